# Applies the "remove warnings, update feature list" edit to open_features.xlsx
# Sheet1 ("20180611") is the active sheet that received content edits:
#   - Row 14: state Open -> on_hold, and a new Details note added
#             ("see en.DM00046982.pdf page 192")
#   - Row 15: Details note updated from "not possible in CortexM family"
#             to "see en.DM00046982.pdf page 192"
#   - Row 17/18: state Open -> on_hold, prio High -> Low
#   - Row 23: prio Low -> High
#   - Row 38: state Open -> on_hold, prio High -> Low
#   - The sheet's active selection/scroll position moved to B38

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlPasteFormats
$xlPasteFormats = -4122

# --- Unmodified cells used purely as format (fill/font) donors -------------
# B34 carries the "no special fill" look (style used by on_hold/Low cells)
# B31 carries the red "High priority" fill
# C2  carries the bold "Details" note font used for D14

# --- New Details note text (inserted first so it gets the lower shared- ---
# --- string index, matching the target 58/59 ordering) ---------------------
$ws.Range("D15").Value = "see en.DM00046982.pdf page 192"
$ws.Range("D14").Value = "see en.DM00046982.pdf page 192"
$ws.Range("C2").Copy()
$ws.Range("D14").PasteSpecial($xlPasteFormats)
$ws.Range("D14").Value = "see en.DM00046982.pdf page 192"

# --- Row 14: State Open -> on_hold ------------------------------------------
$ws.Range("B34").Copy()
$ws.Range("A14").PasteSpecial($xlPasteFormats)
$ws.Range("A14").Value = "on_hold"

# --- Row 17: State Open -> on_hold, Prio High -> Low ------------------------
$ws.Range("B34").Copy()
$ws.Range("A17").PasteSpecial($xlPasteFormats)
$ws.Range("A17").Value = "on_hold"

$ws.Range("B34").Copy()
$ws.Range("B17").PasteSpecial($xlPasteFormats)
$ws.Range("B17").Value = "Low"

# --- Row 18: State Open -> on_hold, Prio High -> Low ------------------------
$ws.Range("B34").Copy()
$ws.Range("A18").PasteSpecial($xlPasteFormats)
$ws.Range("A18").Value = "on_hold"

$ws.Range("B34").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("B18").Value = "Low"

# --- Row 23: Prio Low -> High ------------------------------------------------
$ws.Range("B31").Copy()
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$ws.Range("B23").Value = "High"

# --- Row 38: State Open -> on_hold, Prio High -> Low ------------------------
$ws.Range("B34").Copy()
$ws.Range("A38").PasteSpecial($xlPasteFormats)
$ws.Range("A38").Value = "on_hold"

$ws.Range("B34").Copy()
$ws.Range("B38").PasteSpecial($xlPasteFormats)
$ws.Range("B38").Value = "Low"

# --- Sheet view: scroll to A16, select B38 ----------------------------------
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("B38").Select()
